$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.833.94"
$ws.Range("E2").Value = "  -1.08%  "
$ws.Range("D3").Value = "1.905.93"
$ws.Range("E3").Value = "  -0.36%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.16"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.27%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.003"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.14%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5018"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +3.60%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3817"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.16%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07282"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -1.16%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9070"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -3.04%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.79"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -0.27%  "
$ws.Range("D12").Value = "1.953.73"
$ws.Range("E12").Value = "  -0.01%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07646"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -2.08%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.480"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -0.72%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.607"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -0.51%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "91.38"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -0.47%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.003"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -0.25%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008710"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -1.45%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.003"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.13%  "
$ws.Range("D20").Value = "27.876.99"
$ws.Range("E20").Value = "  -1.01%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.53"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -2.34%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.165"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.09%  "
$ws.Range("E23").Value = "  -0.95%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "154.38"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -1.13%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.870"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -2.45%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.244"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +6.39%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.37"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -1.16%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "115.22"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -1.00%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.906"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -1.26%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.08977"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +0.11%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.204"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -4.12%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.232"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -1.76%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7639"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -1.33%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.642"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -1.00%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02059"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.01%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.548"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -3.26%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.095"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.96%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5585"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +1.50%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.014"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +0.91%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05250"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -1.29%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.950"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -1.17%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.486"
$ws.Range("D42").ClearFormats()
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1512"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -1.18%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "111.55"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +3.47%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.62"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.85%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4797"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.93%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.003"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.21%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.630"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -1.79%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "67.38"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -1.82%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06068"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.73%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.9003"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.12%  "
